$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# The 26/6/2025 row ("J45") used to carry the "no pude adelantar" note.
# That note now belongs to the freshly-reported 27/6/2025 row, so J45
# goes back to "N/A" like the rows above it.
$ws.Range("J45").Value = "N/A"

# Append a new row to the table for 27/6/2025 - the table range grows
# from D4:J45 to D4:J46 automatically.
$newRow = $tbl.ListRows.Add()

# Copy the format of the row above (26/6/2025) down into the new row so
# the new cells inherit the same borders / fonts / number formats.
$ws.Range("D45:J45").Copy()
$ws.Range("D46:J46").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new row's values - same counts as the prior report, new date
# and the note that used to sit on row 45.
$ws.Range("D46").Value = "27/6/2025"
$ws.Range("E46").Value = 297
$ws.Range("F46").Value = 629
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 1012
$ws.Range("J46").Value = "No he podido adelantar, le dare en fin de semana, disculpa (rafael)"

# Re-apply the date format on the new date cell.
$ws.Range("D46").NumberFormat = $ws.Range("D46").NumberFormat
